$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp caption (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 14:05"

# Country-row refresh: rank reshuffles (new leaders inserted, others pushed down)
# and general case-count updates, per the 14:05 data snapshot.
$rowData = @{
    15 = @("India", 71369, 601, 23033, 46026, 0, 16, 2310)
    22 = @("Pakistan", 32674, 1733, 8555, 23395, 111, 57, 724)
    47 = @("Kuwait", 10277, 991, 3101, 7101, 158, 10, 75)
    48 = @("Serbia", 10176, 0, 3290, 6668, 29, 0, 218)
    49 = @("Egipto", 9746, 0, 2172, 7041, 41, 0, 533)
    53 = @("Australia", 6964, 16, 6229, 638, 16, 0, 97)
    61 = @("Ghana", 5127, 427, 494, 4611, 5, 0, 22)
    62 = @("Moldavia", 4995, 0, 1980, 2840, 251, 0, 175)
    63 = @("Afganistan", 4963, 276, 610, 4226, 7, 5, 127)
    76 = @("Guinea", 2213, 67, 771, 1431, 0, 0, 11)
    77 = @("Croacia", 2196, 0, 1784, 321, 13, 0, 91)
    78 = @("Bosnia y Herzegovina", 2158, 17, 1168, 873, 4, 4, 117)
    120 = @("San Marino", 638, 10, 161, 436, 2, 0, 41)
    123 = @("Malta", 506, 3, 434, 67, 1, 0, 5)
    124 = @("Jamaica", 505, 3, 90, 406, 0, 0, 9)
    137 = @("Vietnam", 288, 0, 252, 36, 2, 0, 0)
    215 = @("San Bartolome", 6, 0, 6, 0, 0, 0, 0)
    216 = @("Sahara Occidental", 6, 0, 6, 0, 0, 0, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
